$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.690.59"
$ws.Range("E2").Value = '  +1.57%  '

$ws.Range("D3").Value = "'1.868.38"
$ws.Range("E3").Value = '  +1.80%  '

$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").Value = "'326.71"
$ws.Range("E5").Value = '  -1.15%  '

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").Value = "'0.4628"
$ws.Range("E7").Value = '  +0.63%  '

$ws.Range("D8").Value = "'0.3910"
$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("D9").Value = "'0.07904"
$ws.Range("E9").Value = '  +0.77%  '

$ws.Range("D10").Value = "'0.9716"

$ws.Range("D11").Value = "'22.30"
$ws.Range("E11").Value = '  +2.27%  '

$ws.Range("D12").Value = "'1.912.79"
$ws.Range("E12").Value = '  +3.20%  '

$ws.Range("D13").Value = "'5.735"
$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").Value = "'6.939"
$ws.Range("E14").Value = '  +0.48%  '

$ws.Range("D15").Value = "'0.06922"
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").Value = "'88.59"
$ws.Range("E16").Value = '  +2.39%  '

$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("D18").Value = "'0.00001004"
$ws.Range("E18").Value = '  +1.29%  '

$ws.Range("D19").Value = "'16.94"
$ws.Range("E19").Value = '  +0.51%  '

$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").Value = "'28.663.81"
$ws.Range("E21").Value = '  +1.39%  '

$ws.Range("D22").Value = "'5.325"
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").Value = "'11.08"
$ws.Range("E23").Value = '  +0.84%  '

$ws.Range("D24").Value = "'2.124"
$ws.Range("E24").Value = '  -1.35%  '

$ws.Range("D25").Value = "'2.113.57"
$ws.Range("E25").Value = '  +2.83%  '

$ws.Range("D26").Value = "'155.21"
$ws.Range("E26").Value = '  +1.44%  '

$ws.Range("D27").Value = "'19.33"
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").Value = "'5.777"
$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").Value = "'1.996"
$ws.Range("E29").Value = '  +1.31%  '

$ws.Range("D30").Value = "'119.24"
$ws.Range("E30").Value = '  +2.21%  '

$ws.Range("D31").Value = "'0.09353"
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").Value = "'0.9381"
$ws.Range("E32").Value = '  -0.27%  '

$ws.Range("D33").Value = "'5.324"
$ws.Range("E33").Value = '  +0.81%  '

$ws.Range("D34").Value = "'1.339"
$ws.Range("E34").Value = '  +1.04%  '

$ws.Range("E35").Value = '  -2.89%  '

$ws.Range("D36").Value = "'0.05827"
$ws.Range("E36").Value = '  -3.43%  '

$ws.Range("D37").Value = "'0.02115"
$ws.Range("E37").Value = '  -2.08%  '

$ws.Range("D38").Value = "'1.156"
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("D39").Value = "'7.886"
$ws.Range("E39").Value = '  +4.58%  '

$ws.Range("D40").Value = "'0.5665"
$ws.Range("E40").Value = '  +0.99%  '

$ws.Range("D41").Value = "'9.947"
$ws.Range("E41").Value = '  -0.49%  '

$ws.Range("D42").Value = "'0.1779"
$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'2.237"
$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = "'0.07259"
$ws.Range("E44").Value = '  +3.36%  '

$ws.Range("D45").Value = "'11.72"
$ws.Range("E45").Value = '  +0.95%  '

$ws.Range("D46").Value = "'0.5323"
$ws.Range("E46").Value = '  +0.49%  '

$ws.Range("E47").Value = '  -8.34%  '

$ws.Range("D48").Value = "'1.850"
$ws.Range("E48").Value = '  +0.55%  '

$ws.Range("D49").Value = "'113.60"
$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("D50").Value = "'2.351"
$ws.Range("E50").Value = '  +1.55%  '

$ws.Range("D51").Value = "'1.005"
$ws.Range("E51").Value = '  +0.33%  '
